$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45: InceptionResNetV2_19
$ws.Cells.Item(45, 1).Value = "InceptionResNetV2_19"
$ws.Cells.Item(45, 2).Value = "InceptionResNetV2"
$ws.Cells.Item(45, 3).Value = 0.95711892797319931
$ws.Cells.Item(45, 4).Value = 0.91600000000000004
$ws.Cells.Item(45, 5).Value = 7
$ws.Cells.Item(45, 6).Value = "['Functional', 'Dropout', 'BatchNormalization', 'Dense', 'Dropout', 'BatchNormalization', 'Dense']"
$ws.Cells.Item(45, 7).Value = "[0.1, 0.1]"
$ws.Cells.Item(45, 8).Value = "[('relu', 1024), ('softmax', 15)]"
$ws.Cells.Item(45, 9).Value = 100
$ws.Cells.Item(45, 10).Value = 0.00001
$ws.Cells.Item(45, 11).Value = "{'monitor': 'val_loss', 'patience': 30, 'min_delta': 0, 'restore_best_weights': True}"
$ws.Cells.Item(45, 12).Value = "{'Train': 32, 'Validation': 32}"
$ws.Cells.Item(45, 13).Value = "{'zoom_range': 0.2, 'rotation_range': 30, 'shear_range': 0.2, 'brightness_range': None, 'horizontal_flip': True, 'width_shift_range': 0.2, 'height_shift_range': 0.2}"
$ws.Cells.Item(45, 14).Value = 55936239
$ws.Cells.Item(45, 15).Value = 0
$ws.Cells.Item(45, 16).Value = 6391.827333688736

# Row 46: InceptionResNetV2_20
$ws.Cells.Item(46, 1).Value = "InceptionResNetV2_20"
$ws.Cells.Item(46, 2).Value = "InceptionResNetV2"
$ws.Cells.Item(46, 3).Value = 0.97185929648241209
$ws.Cells.Item(46, 4).Value = 0.92333333333333334
$ws.Cells.Item(46, 5).Value = 7
$ws.Cells.Item(46, 6).Value = "['Functional', 'Dropout', 'BatchNormalization', 'Dense', 'Dropout', 'BatchNormalization', 'Dense']"
$ws.Cells.Item(46, 7).Value = "[0.1, 0.1]"
$ws.Cells.Item(46, 8).Value = "[('relu', 1024), ('softmax', 15)]"
$ws.Cells.Item(46, 9).Value = 200
$ws.Cells.Item(46, 10).Value = 0.00001
$ws.Cells.Item(46, 11).Value = "{'monitor': 'val_loss', 'patience': 25, 'min_delta': 0, 'restore_best_weights': True}"
$ws.Cells.Item(46, 12).Value = "{'Train': 32, 'Validation': 32}"
$ws.Cells.Item(46, 13).Value = "{'zoom_range': 0.2, 'rotation_range': 30, 'shear_range': 0.2, 'brightness_range': None, 'horizontal_flip': True, 'width_shift_range': 0.2, 'height_shift_range': 0.2}"
$ws.Cells.Item(46, 14).Value = 55936239
$ws.Cells.Item(46, 15).Value = 0
$ws.Cells.Item(46, 16).Value = 8905.7663018703461

# Widen column A to fit the longer model names, and scroll the view down
# to roughly where the newly added rows are.
$ws.Columns.Item(1).ColumnWidth = 45.38
$ws.Range("A25").Select()
